$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.270.42'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.685.33'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.26'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5252'
$ws.Range('E6').Value = '  +2.94%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2715'
$ws.Range('E8').Value = '  +2.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06420'
$ws.Range('E9').Value = '  +0.55%  '
$ws.Range('E10').Value = '  +2.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07493'
$ws.Range('E11').Value = '  +1.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.701.54'
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.557'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5801'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008466'
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.28'
$ws.Range('E16').Value = '  -0.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.331.58'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.926'
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.87'
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.08'
$ws.Range('E21').Value = '  -0.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.194'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '144.54'
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.708'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1237'
$ws.Range('E26').Value = '  +4.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.79'
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06637'
$ws.Range('E28').Value = '  +12.33%  '
$ws.Range('E29').Value = '  +5.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.328'
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.575'
$ws.Range('E31').Value = '  +1.70%  '
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.662'
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.028'
$ws.Range('E34').Value = '  +1.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6200'
$ws.Range('E35').Value = '  +3.08%  '
$ws.Range('E36').Value = '  +1.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.704'
$ws.Range('E37').Value = '  +1.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.388'
$ws.Range('E38').Value = '  +5.73%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01621'
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.106.85'
$ws.Range('E40').Value = '  +3.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8758'
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('E42').Value = '  +0.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.69'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.832.83'
$ws.Range('E44').Value = '  +0.71%  '
$ws.Range('E45').Value = '  -2.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.74'
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.189'
$ws.Range('E47').Value = '  +1.71%  '
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('E49').Value = '  +1.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4304'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.035'
$ws.Range('E51').Value = '  +2.60%  '
